$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 204.6
$ws.Range("C2").Value = 409.2

$ws.Range("B3").Value = 102.3
$ws.Range("C3").Value = 204.6

$ws.Range("B4").Value = 51.15
$ws.Range("C4").Value = 102.3

$ws.Range("B5").Value = 153.45
$ws.Range("C5").Value = 306.9
